$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of a subset of columns (D, M, N, O, P, R, S)
# across rows 2, 3 and 5: row2 <- old row3, row3 <- old row5, row5 <- old row2.
# Capture the "before" values first (using Value2 to get raw data, not a
# COM wrapper description), then write them back in rotated order.

$row2_D = $ws.Range("D2").Value2
$row2_M = $ws.Range("M2").Value2
$row2_N = $ws.Range("N2").Value2
$row2_O = $ws.Range("O2").Value2
$row2_P = $ws.Range("P2").Value2
$row2_R = $ws.Range("R2").Value2
$row2_S = $ws.Range("S2").Value2

$row3_D = $ws.Range("D3").Value2
$row3_M = $ws.Range("M3").Value2
$row3_N = $ws.Range("N3").Value2
$row3_O = $ws.Range("O3").Value2
$row3_P = $ws.Range("P3").Value2
$row3_R = $ws.Range("R3").Value2
$row3_S = $ws.Range("S3").Value2

$row5_D = $ws.Range("D5").Value2
$row5_M = $ws.Range("M5").Value2
$row5_N = $ws.Range("N5").Value2
$row5_O = $ws.Range("O5").Value2
$row5_P = $ws.Range("P5").Value2
$row5_R = $ws.Range("R5").Value2
$row5_S = $ws.Range("S5").Value2

# Row 2 gets old row 3 values
$ws.Range("D2").Value = $row3_D
$ws.Range("M2").Value = $row3_M
$ws.Range("N2").Value = $row3_N
$ws.Range("O2").Value = $row3_O
$ws.Range("P2").Value = $row3_P
$ws.Range("R2").Value = $row3_R
$ws.Range("S2").Value = $row3_S

# Row 3 gets old row 5 values
$ws.Range("D3").Value = $row5_D
$ws.Range("M3").Value = $row5_M
$ws.Range("N3").Value = $row5_N
$ws.Range("O3").Value = $row5_O
$ws.Range("P3").Value = $row5_P
$ws.Range("R3").Value = $row5_R
$ws.Range("S3").Value = $row5_S

# Row 5 gets old row 2 values
$ws.Range("D5").Value = $row2_D
$ws.Range("M5").Value = $row2_M
$ws.Range("N5").Value = $row2_N
$ws.Range("O5").Value = $row2_O
$ws.Range("P5").Value = $row2_P
$ws.Range("R5").Value = $row2_R
$ws.Range("S5").Value = $row2_S
